# Updates cryptos list values (Price and Volume(1h) columns) to reflect
# the latest scrape, plus refreshes row 51 coin (dogwifhat -> WhiteBITCoin).

function Set-TextValue($sheet, $ref, $val) {
    # Force the cell to be written as text even when the value looks
    # like a plain number (e.g. "580.33"), then drop the temporary
    # number-format so the cell keeps its original (default) style.
    $sheet.Range($ref).NumberFormat = "@"
    $sheet.Range($ref).Value = $val
    $sheet.Range($ref).ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.098.31"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "2.555.87"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws "D5" "580.33"
$ws.Range("E5").Value = "  +1.77%  "
Set-TextValue $ws "D6" "147.28"
$ws.Range("E7").Value = "  +0.00%  "
Set-TextValue $ws "D8" "0.585"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("E11").Value = "  -0.03%  "
Set-TextValue $ws "D12" "0.354"
$ws.Range("E12").Value = "  -0.92%  "
Set-TextValue $ws "D13" "27.55"
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("D14").Value = "3.011.51"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").Value = "63.046.11"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "2.542.94"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("E18").Value = "  -2.54%  "
Set-TextValue $ws "D19" "339.37"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("E20").Value = "  +0.01%  "
Set-TextValue $ws "D21" "6.77"
$ws.Range("E21").Value = "  -0.35%  "
Set-TextValue $ws "D22" "0.999"
$ws.Range("E22").Value = "  -0.11%  "
Set-TextValue $ws "D23" "65.58"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").Value = "2.675.52"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("E27").Value = "  -4.48%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("E29").Value = "  -1.12%  "
Set-TextValue $ws "D30" "7.70"
$ws.Range("E30").Value = "  +7.25%  "
$ws.Range("E31").Value = "  +4.76%  "
$ws.Range("E32").Value = "  -0.18%  "
Set-TextValue $ws "D33" "178.05"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  -0.84%  "
Set-TextValue $ws "D35" "420.76"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("E36").Value = "  -0.87%  "
Set-TextValue $ws "D37" "19.15"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  +0.48%  "
Set-TextValue $ws "D43" "150.93"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("E44").Value = "  -0.03%  "
Set-TextValue $ws "D45" "20.80"
$ws.Range("E45").Value = "  +0.02%  "
Set-TextValue $ws "D46" "0.0540"
$ws.Range("E46").Value = "  +3.05%  "
$ws.Range("E47").Value = "  -0.95%  "
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").Value = "  +0.03%  "
Set-TextValue $ws "D50" "18.30"
$ws.Range("E50").Value = "  -1.24%  "

# Row 51: coin swapped out for a new entry
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws "D51" "11.37"
$ws.Range("E51").Value = "  -0.26%  "
